$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Policy")

# Rename "Policy" -> "policy"
$ws.Name = "policy"

# Insert a new first column (old A -> B, old B -> C) and add the new
# "table elements" lookup column down the left side of the sheet.
$ws.Columns.Item(1).Insert()

$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = "table elements"
$ws.Range("A1").Interior.Color = 65535

$ws.Range("A2").Value = "POLICY ID"
$ws.Range("A3").Value = "TERM"
$ws.Range("A4").Value = "TOTAL AMOUNT"
$ws.Range("A5").Value = "PER MONTH"
$ws.Range("A6").Value = "PAYMENT METHOD"
$ws.Range("A7").Value = "COVERAGE"
$ws.Range("A8").Value = "AGE LIMIT"

# The shifted-over cells (old col A values) keep their style; the brand
# new label cells in rows 2-4 go back to the sheet's default style.
$ws.Range("A2:A4").Style = "Normal"

$ws.Columns.Item(1).ColumnWidth = 19.5

# Make "policy" the active sheet / tab, with I21 selected - this also
# clears tabSelected on whichever sheet previously had it (Agents).
$ws.Activate() | Out-Null
$ws.Range("I21").Select() | Out-Null
